$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate: wipe every existing value on the sheet
# (this removes the old header/row2/row6 content but keeps any
# pre-existing cell formatting such as the Hyperlink style on A2).
$ws.Cells.ClearContents()

# ---- Row 1 : headers -----------------------------------------------------
$ws.Range("A1").Value = "link"
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "CM"
$ws.Range("D1").Value = "title"
$ws.Range("E1").Value = "department"
$ws.Range("F1").Value = "location"
$ws.Range("G1").Value = "deadline"
$ws.Range("H1").Value = "number"
$ws.Range("I1").Value = "interlinkregex"
$ws.Range("J1").Value = "finallinkregex"

# ---- Row 2 : cmbchina (campus) -------------------------------------------
$ws.Range("A2").Value = "http://career.cmbchina.com/Campus/Campus.aspx"
$ws.Range("B2").Value = "cmbchina"
$ws.Range("C2").Value = "C"
$ws.Range("D2").Value = '//*[@id="rightdiv"]/div/div[1]/text()'
$ws.Range("E2").Value = '//*[@id="rightdiv"]/div/div[2]/p[3]/text()'
$ws.Range("F2").Value = '//*[@id="rightdiv"]/div/div[2]/p[5]/text()'
$ws.Range("G2").Value = '//*[@id="rightdiv"]/div/div[2]/p[6]/text()'
$ws.Range("I2").Value = "branch="
$ws.Range("J2").Value = "Position.aspx.id"

# ---- Row 3 : cmbchina (social) -------------------------------------------
$ws.Range("A3").Value = "http://career.cmbchina.com/Social/Default.aspx"
$ws.Range("B3").Value = "cmbchina"
$ws.Range("C3").Value = "M"
$ws.Range("D3").Value = '//*[@id="rightdiv"]/div[1]/div[1]/text()'
$ws.Range("E3").Value = '//*[@id="rightdiv"]/div[1]/div[2]/p[3]/text()'
$ws.Range("F3").Value = '//*[@id="rightdiv"]/div[1]/div[2]/p[5]/text()'
$ws.Range("G3").Value = '//*[@id="rightdiv"]/div[1]/div[2]/p[6]/text()'
$ws.Range("I3").Value = "branch="
$ws.Range("J3").Value = "Position.aspx.id"

# ---- Row 4 : ccb -----------------------------------------------------------
$ws.Range("A4").Value = "http://job.ccb.com/ccbjob/cn/job/index.jsp"
$ws.Range("B4").Value = "ccb"
$ws.Range("C4").Value = "C"
$ws.Range("D4").Value = '//*[@id="title"]/strong/text()'
$ws.Range("E4").Value = '//*[@id="title"]/strong/text()'
$ws.Range("F4").Value = '//*[@id="data"]/table/tbody/tr[2]/td/p[5]/font/text()'
$ws.Range("I4").Value = "branch_notice_list2"
$ws.Range("J4").Value = "info"

# ---- Row 5 : abchina --------------------------------------------------------
$ws.Range("B5").Value = "abchina"
$ws.Range("C5").Value = "C"
$ws.Range("D5").Value = '//*[@id="frmPutQuestion"]/table[3]/tr[2]/td[2]/text()'
$ws.Range("E5").Value = '//*[@id="frmPutQuestion"]/table[3]/tr[6]/td[2]/text()'
$ws.Range("F5").Value = '//*[@id="frmPutQuestion"]/table[3]/tbody/tr[3]/td[2]/text()'
$ws.Range("G5").Value = '//*[@id="frmPutQuestion"]/table[3]/tr[5]/td[2]/text()'
$ws.Range("I5").Value = "jobPublish"
$ws.Range("J5").Value = "jobDetails"

# A5 carries a real hyperlink (like A2), so add it via the Hyperlinks
# collection and then restore the built-in "Hyperlink" cell style
# (Hyperlinks.Add on its own allocates a fresh, slightly different xf).
$ws.Hyperlinks.Add($ws.Range("A5"), "http://job.abchina.com/rio/index.do?action=openHome&toId=2") | Out-Null
$ws.Range("A5").Style = "Hyperlink"

# Match the saved selection state from the workbook.
[void]$ws.Range("J5").Select()
